# Update the pricing workbook for "xuất báo giá PDF":
#  1. Rename the first sheet tab, removing the stray space before "(CTN009)".
#  2. Fix the "Có VAT" / "Chưa VAT" header labels so the first (base-price)
#     table reads "Chưa VAT" and the second (computed, *1.08) table reads
#     "Có VAT" on the "Tiêu chuẩn TMĐT ĐG(CTN009)" sheet (it previously had
#     "Có VAT" duplicated on both tables).
#  3. Reset the active selection on that sheet back to A1.

$wb = $excel.ActiveWorkbook

$wsDG = $wb.Worksheets.Item(1)
$wsCTN007 = $wb.Worksheets.Item(2)

# 1. Rename sheet tab (drop the space before the opening parenthesis).
$wsDG.Name = "Tiêu chuẩn TMĐT ĐG(CTN009)"

$coVat = "Nấc cước" + [char]10 + "(Có VAT)"
$chuaVat = "Nấc cước" + [char]10 + "(Chưa VAT)"

# 2. Correct the VAT labels.
# "Tiêu chuẩn TMĐT ĐG(CTN009)" sheet: base table (row 1) -> Chưa VAT,
# computed table (row 7, values = ROUNDUP(base*1.08)) -> Có VAT.
$wsDG.Range("A1").Value = $chuaVat
$wsDG.Range("A7").Value = $coVat

# "Tiêu chuẩn TMĐT(CTN007)" sheet already has the correct labels, but set
# them explicitly so both tables stay in sync.
$wsCTN007.Range("A1").Value = $chuaVat
$wsCTN007.Range("A8").Value = $coVat

# 3. Reset the selection on the renamed sheet to A1.
$wsDG.Activate()
$wsDG.Range("A1").Select()
